$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: new blank-dates row, only the duration formula (D32-C32) present.
$ws.Range("E32").Formula = "=D32-C32"
$ws.Range("E32").NumberFormat = "h:mm:ss"

# Row 33 already has B33 (date 2017-07-17). Add the clock-in/out pair + duration.
$ws.Range("C33").Value = 0.47916666666666669
$ws.Range("C33").NumberFormat = "h:mm"
$ws.Range("D33").Value = 0.55208333333333337
$ws.Range("D33").NumberFormat = "h:mm"
$ws.Range("E33").Formula = "=D33-C33"
$ws.Range("E33").NumberFormat = "h:mm:ss"

# Row 34: second clock-in/out pair for the same day + duration.
$ws.Range("C34").Value = 0.57638888888888895
$ws.Range("C34").NumberFormat = "h:mm"
$ws.Range("D34").Value = 0.72222222222222221
$ws.Range("D34").NumberFormat = "h:mm"
$ws.Range("E34").Formula = "=D34-C34"
$ws.Range("E34").NumberFormat = "h:mm:ss"

# Row 35: trailing blank row that still belongs to the shared-formula block,
# so it keeps the duration column's styling but no value/formula.
$ws.Range("E35").Formula = "=D35-C35"
$ws.Range("E35").ClearContents()
$ws.Range("E35").NumberFormat = "h:mm:ss"

# Move the active selection to match the saved view state.
$ws.Range("G34").Select() | Out-Null
